$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137 (weekly price update for "Hortaliza,
# Terminal La Palmera de La Serena - Espinaca"), shifting the existing
# rows 137-175 down to 138-176 and preserving their data/format.
$ws.Rows(137).Insert()

# Populate the newly inserted row 137 with the new weekly record.
$ws.Cells.Item(137, 1).Value = 8
$ws.Cells.Item(137, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(137, 3).Value = "Coquimbo"
$ws.Cells.Item(137, 4).Value = 44508
$ws.Cells.Item(137, 5).Value = 4
$ws.Cells.Item(137, 6).Value = 100112012
$ws.Cells.Item(137, 7).Value = "Espinaca"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 400
$ws.Cells.Item(137, 12).Value = 500
$ws.Cells.Item(137, 13).Value = 450
$ws.Cells.Item(137, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(137, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(137, 16).Value = 900
$ws.Cells.Item(137, 17).Value = 0.5
$ws.Cells.Item(137, 18).Value = "Hortaliza"
